# "hours update and TAR update"
# Append two new status-report rows (18 and 19) to Sheet1, both dated
# 1/26/2010, mirroring the existing "Date" column entries that are stored
# as plain text (e.g. row 17's "1/21/2010") rather than true Excel dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use a leading apostrophe so Excel stores the date as text (matching the
# existing text-style date entries in column A) instead of auto-converting
# it to a date serial number.
$ws.Range("A18").Value = "'1/26/2010"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Group Meeting"

$ws.Range("A19").Value = "'1/26/2010"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "Weekly Meeting"

# Drop the quote-prefix formatting the apostrophe entry introduced so the
# new cells keep the sheet's default (unstyled) look, same as row 17.
$ws.Range("A18:A19").ClearFormats()

# Move the active selection down to A20, just past the newly added rows.
$ws.Range("A20").Select()
